$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(175)
$r = $p.Range
$r.InsertParagraphAfter()
for ($i = 172; $i -le 184; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $rr = $pp.Range
    Write-Output "Index: $i Start: $($rr.Start) End: $($rr.End) Text:[$($rr.Text)]"
}
